# Add grade names, both valid and invalid, to the DataSource sheet's
# "Grade" column (column E), replacing three existing grade codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSource")

# Row 11: was "03" -> "Third"
$ws.Range("E11").Value = "Third"

# Row 19: was "02" -> "Grade 2"
$ws.Range("E19").Value = "Grade 2"

# Row 26: was "K" -> "Grade 20"
$ws.Range("E26").Value = "Grade 20"
